$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36: GFG - Find the Closest Element in BST
$ws.Range("A36").Value = "NA"
$ws.Range("A36").HorizontalAlignment = -4131

$ws.Range("B36").Value = "GFG"
$ws.Range("B36").HorizontalAlignment = -4131
$ws.Range("B36").VerticalAlignment = -4160
$ws.Range("B36").WrapText = $true

$ws.Range("C36").Value = "Find the Closest Element in BST"
$ws.Range("C36").HorizontalAlignment = -4131

# Row 37: Floor and Ceil in BST.py
$ws.Range("A37").Value = "NA"
$ws.Range("A37").HorizontalAlignment = -4131

$ws.Range("B37").Value = "NA"
$ws.Range("B37").HorizontalAlignment = -4131
$ws.Range("B37").VerticalAlignment = -4160
$ws.Range("B37").WrapText = $true

$ws.Range("C37").Value = "Floor  and Ceil in BST.py"

# Update the active cell selection to reflect the new last row
$ws.Range("A38").Select()
